# The presentation currently uses the "Integral" theme (green colour
# scheme) for its slide master / overall design. The target edit swaps
# the deck back to the stock "Office Theme" colour scheme.
#
# PowerPoint's theme colours are addressed positionally through the
# MsoThemeColorSchemeIndex ordering:
#   1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
#   9=accent5 10=accent6 11=hlink 12=folHlink
# `.RGB` uses the VBA packed integer form (R + G*256 + B*65536), so a
# small helper converts plain hex colours into that form.

function ConvertTo-RgbValue([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

function ConvertHex-ToRgbValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ConvertTo-RgbValue $r $g $b
}

# Stock "Office Theme" colour scheme values.
$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

$colorCount = $themeColors.Count
if (-not $colorCount) {
    $colorCount = $officeThemeColors.Count
}

for ($i = 1; $i -le $colorCount; $i++) {
    $themeColors.Colors($i).RGB = ConvertHex-ToRgbValue $officeThemeColors[$i - 1]
}
